# Apply data corrections to the "R11_Y" (column R) and "R11_Z" (column W)
# statistics on the "Sheet 1" worksheet, per updated ESS11 edition 4.0 data.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet 1")

$ws.Range("R3").Value = 0.299
$ws.Range("W3").Value = 0.105
$ws.Range("R4").Value = 2.514
$ws.Range("R5").Value = 4.45
$ws.Range("R6").Value = 649.244
$ws.Range("W6").Value = 0.323
$ws.Range("R7").Value = 30.239
$ws.Range("W7").Value = 0.792
$ws.Range("R8").Value = 33.991
$ws.Range("W8").Value = 0.848
$ws.Range("R9").Value = 12.411
$ws.Range("W9").Value = 0.317
$ws.Range("R10").Value = 1952.278
$ws.Range("R11").Value = 1.908
$ws.Range("W11").Value = 0.989
$ws.Range("R12").Value = 4.983
$ws.Range("W12").Value = 0.926
$ws.Range("R13").Value = 2.089
$ws.Range("W13").Value = 0.973
$ws.Range("R14").Value = 2.057
$ws.Range("W14").Value = 0.981
$ws.Range("R15").Value = 7.268
$ws.Range("W15").Value = 0.995
$ws.Range("R16").Value = 2.186
$ws.Range("R17").Value = 5.143
$ws.Range("W17").Value = 0.967
$ws.Range("R18").Value = 5.162
$ws.Range("W18").Value = 0.969
$ws.Range("R19").Value = 4.825
$ws.Range("W19").Value = 0.965
$ws.Range("R20").Value = 3.849
$ws.Range("W20").Value = 0.888
$ws.Range("R21").Value = 4.344
$ws.Range("W21").Value = 0.853
$ws.Range("R22").Value = 2.654
$ws.Range("R23").Value = 5.53
$ws.Range("W23").Value = 0.993
$ws.Range("R24").Value = 4.936
$ws.Range("R25").Value = 4.948
$ws.Range("R26").Value = 4.636
$ws.Range("W26").Value = 0.991
$ws.Range("R27").Value = 4.763
$ws.Range("W27").Value = 0.966
$ws.Range("R28").Value = 4.261
$ws.Range("W28").Value = 0.982
$ws.Range("R29").Value = 4.962
$ws.Range("W29").Value = 0.955
$ws.Range("R30").Value = 3.842
$ws.Range("R31").Value = 5.022
$ws.Range("W31").Value = 0.991
$ws.Range("R32").Value = 6.935
$ws.Range("R33").Value = 4.222
$ws.Range("W33").Value = 0.935
$ws.Range("R34").Value = 5.132
$ws.Range("W34").Value = 0.977
$ws.Range("R35").Value = 6.183
$ws.Range("W35").Value = 0.988
$ws.Range("R36").Value = 3.435
$ws.Range("W36").Value = 0.984
$ws.Range("R37").Value = 4.214
$ws.Range("R38").Value = 3.409
$ws.Range("W38").Value = 0.981
$ws.Range("R39").Value = 4.54
$ws.Range("W39").Value = 0.926
$ws.Range("R40").Value = 5.451
$ws.Range("W40").Value = 0.891
$ws.Range("R41").Value = 0.254
$ws.Range("W41").Value = 0.992
$ws.Range("R42").Value = 0.909
$ws.Range("W42").Value = 0.999
$ws.Range("R43").Value = 0.353
$ws.Range("W43").Value = 0.668
$ws.Range("R44").Value = 0.421
$ws.Range("W44").Value = 0.973
$ws.Range("R46").Value = 0.029
$ws.Range("W46").Value = 0.464
$ws.Range("R51").Value = 0.011
$ws.Range("R52").Value = 0.014
$ws.Range("R53").Value = 0.084
$ws.Range("W53").Value = 0.989
$ws.Range("R54").Value = 0.009
$ws.Range("R56").Value = 0.914
$ws.Range("R57").Value = 0.018
$ws.Range("R58").Value = 0.015
$ws.Range("R59").Value = 0.011
$ws.Range("R63").Value = 0.163
$ws.Range("R64").Value = 0.771
$ws.Range("W64").Value = 0.895
$ws.Range("R66").Value = 0.459
$ws.Range("W66").Value = 0.999
$ws.Range("R67").Value = 0.297
$ws.Range("W67").Value = 0.866
$ws.Range("R69").Value = 0.341
$ws.Range("W69").Value = 0.435
$ws.Range("R70").Value = 0.633
$ws.Range("R71").Value = 0.104
$ws.Range("W71").Value = 0.361
$ws.Range("R73").Value = 0.11
$ws.Range("W73").Value = 0.265
$ws.Range("R74").Value = 0.267
$ws.Range("R75").Value = 0.102
$ws.Range("W75").Value = 0.265
$ws.Range("R76").Value = 0.725
$ws.Range("R77").Value = 0.063
$ws.Range("W77").Value = 0.898
